$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 329; this shifts the existing rows
# 329:482 down to 330:483 (matching the target dimension A1:R483).
$ws.Rows.Item(329).Insert()

# Populate the newly inserted row 329 with the new weekly record.
$ws.Cells.Item(329, 1).Value = 8
$ws.Cells.Item(329, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(329, 3).Value = "Coquimbo"
$ws.Cells.Item(329, 4).Value = 45205
$ws.Cells.Item(329, 5).Value = 4
$ws.Cells.Item(329, 6).Value = 100112012
$ws.Cells.Item(329, 7).Value = "Espinaca"
$ws.Cells.Item(329, 8).Value = "Sin especificar"
$ws.Cells.Item(329, 9).Value = "Primera"
$ws.Cells.Item(329, 10).Value = 1480
$ws.Cells.Item(329, 11).Value = 450
$ws.Cells.Item(329, 12).Value = 500
$ws.Cells.Item(329, 13).Value = 475
$ws.Cells.Item(329, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(329, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(329, 16).Value = 950
$ws.Cells.Item(329, 17).Value = 0.5
$ws.Cells.Item(329, 18).Value = "Hortaliza"
